# Update the percentage metrics sheet:
#  1. Insert a new model row ("llama3_70b_instruct") in its alphabetical/
#     logical position (right after "llama3_2_3b_instruct"), shifting the
#     following rows down - robust to the row already being present.
#  2. Re-round every percentage value in the data range to 2 decimal
#     places (the metrics-generation script now writes rounded numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data to make sure is present in the sheet.
$newModelName = "llama3_70b_instruct"
$newModelValues = @(68.01, 71.89, 79.23, 64.36, 83.84, 88.55, 95.27)

# Find the last used row/column so this keeps working even if the sheet
# already has extra rows (robust to existing files).
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# Locate the row for $newModelName if it already exists, and figure out
# where it should be inserted (first row, after the header, whose model
# name sorts after $newModelName) if it doesn't.
$existingRow = 0
$insertBeforeRow = 0
for ($r = 2; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    if ($name -eq $newModelName) {
        $existingRow = $r
    }
    if ($insertBeforeRow -eq 0 -and $name -gt $newModelName) {
        $insertBeforeRow = $r
    }
}

if ($existingRow -eq 0) {
    if ($insertBeforeRow -eq 0) {
        $insertBeforeRow = $lastRow + 1
    }
    $ws.Rows.Item($insertBeforeRow).Insert()
    $targetRow = $insertBeforeRow
    $lastRow = $lastRow + 1
} else {
    $targetRow = $existingRow
}

$ws.Cells.Item($targetRow, 1).Value2 = $newModelName
for ($i = 0; $i -lt $newModelValues.Length; $i++) {
    $ws.Cells.Item($targetRow, 2 + $i).Value2 = $newModelValues[$i]
}

# Refresh used range extent after the possible insert.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# Round every numeric metric cell (everything but the header row / name
# column) to 2 decimal places.
for ($r = 2; $r -le $lastRow; $r++) {
    for ($c = 2; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -ne $null) {
            $cell.Value2 = [Math]::Round([double]$v, 2)
        }
    }
}
